# Stats.xlsx update — "Add files via upload"
# Adds a new entrant ("heccinhayley") plus the knock-on count bumps that
# come from re-tallying the BR stats after the addition, and records the
# organizer credit for "grandpaszabo" on row 36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Season summary (A/B columns): S8 entrant count goes 5 -> 6, and the
#     SUM(B2:B12) total recalculates to 94 automatically, but we set the
#     cached value explicitly to match too.
$ws.Range("B10").Value = 6

# --- Organizers list (H/I columns): new credit for grandpaszabo on row 36
$ws.Range("H36").Value = "grandpaszabo"
$ws.Range("I36").Value = 1

# --- All Entrants BR counts (D/E columns): each of these existing
#     entrants' tallies increments by 1
$ws.Range("E40").Value = 2    # cokes
$ws.Range("E41").Value = 6    # crimsonavix
$ws.Range("E42").Value = 2    # crump79
$ws.Range("E43").Value = 7    # custom
$ws.Range("E45").Value = 31   # malkier
$ws.Range("E48").Value = 9    # dem gaming
$ws.Range("E70").Value = 4    # grandpaszabo
$ws.Range("E80").Value = 31   # jessandy
$ws.Range("E81").Value = 11   # Jkoper
$ws.Range("E87").Value = 57   # Kingdahl
$ws.Range("E100").Value = 8   # moneymerks
$ws.Range("E121").Value = 2   # pool float
$ws.Range("E122").Value = 39  # potatopony
$ws.Range("E141").Value = 12  # specialk3782
$ws.Range("E163").Value = 9   # ultralavos
$ws.Range("E166").Value = 3   # von

# --- New entrant row appended to the All Entrants list
$ws.Range("D175").Value = "heccinhayley"
$ws.Range("E175").Value = 1

# --- Recalculate formulas (B13 = SUM(B2:B12)) and restore the saved
#     selection/viewport from the edit session
$wb.Application.Calculate() | Out-Null
$ws.Range("F43").Select() | Out-Null
